$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# URLs of the existing references, keyed by their original row number
# (rows 2..16 before the new entry is inserted).
$urls = @{
    2  = "https://ideas.repec.org/p/pra/mprapa/79809.html"
    3  = "https://ideas.repec.org/p/pra/mprapa/83154.html"
    4  = "https://doi.org/10.1080/02692171.2019.1645816"
    5  = "https://www.bcentral.cl/en/content/-/details/monetary-policy-report-june-2015"
    6  = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4043"
    7  = "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/documento-de-trabajo-n-883"
    8  = "https://ideas.repec.org/a/chb/bcchni/v15y2012i1p105-117.html"
    9  = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4042"
    10 = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4041"
    11 = "https://ideas.repec.org/p/chb/bcchep/56.html"
    12 = "https://ideas.repec.org/a/cml/moneta/vxxxiiy2009i2p181-208.html"
    13 = "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwiE7Y60i5TvAhVQErkGHfecC4MQFjABegQIAhAD&url=http%3A%2F%2Fwww.oecd.org%2Fstd%2Fleading-indicators%2F43815334.pdf&usg=AOvVaw3BstLuhLtAOtjJeL5SsMj4"
    14 = "https://www.sciencedirect.com/science/article/abs/pii/S0169207019300676"
    15 = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/3564"
    16 = "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4867"
}
$newUrl = "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/analisis-de-sentimiento-basado-en-el-informe-de-percepciones-de-negocios-del-banco-central-de-chile"

# Remove the existing hyperlink objects; we will rebuild them in the final
# row order once the new row has been inserted.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Hyperlinks.Delete()
}

# --- Insert the new literature entry as row 10 ---------------------------
# ("del Pilar M / Peralta H / Avila B", 2020) belongs right after
# "Chanut N" and before "Figueroa C" when the table is kept sorted by
# Author 1 then Year, so it is inserted at row 10 and everything that used
# to be row 10 onward shifts down by one.
$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = "del Pilar M"
$ws.Range("B10").Value = "Peralta H"
$ws.Range("C10").Value = "Ávila B"
$ws.Range("E10").Value = 2020
$ws.Range("E10").Style = $ws.Range("E11").Style
$ws.Range("F10").Value = "Análisis de Sentimiento Basado en el Informe de Percepciones de Negocios del Banco Central de Chile"
$ws.Range("G10").Value = "Working Paper"
$ws.Range("H10").Value = "Documentos de Trabajo (Banco Central)"
$ws.Range("I10").Value = "Macroeconomía"
$ws.Range("J10").Value = $newUrl

# --- Re-create all the hyperlinks in the final row order ------------------
for ($r = 2; $r -le 9; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 10), $urls[$r])
    $ws.Cells.Item($r, 10).Style = "Hipervínculo"
}
for ($r = 11; $r -le 17; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 10), $urls[$r - 1])
    $ws.Cells.Item($r, 10).Style = "Hipervínculo"
}
$ws.Hyperlinks.Add($ws.Range("J10"), $newUrl)
$ws.Range("J10").Style = "Hipervínculo"

$ws.Range("B10").Select()
